$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row before row 89 by duplicating it (keeps styles/number
# formats intact), which pushes the old row 89 (and everything below it)
# down by one row.
$ws.Rows.Item(89).Copy()
$ws.Rows.Item(89).Insert()

# Overwrite the new row 89 with this week's new record. The columns that
# are not listed here (A, B, C, E, F, G, H, I, N, Q, R) keep the values
# copied from the row that used to be row 89, since the diff shows them
# unchanged.
$ws.Cells.Item(89, 4).Value = 44546          # D89 Fecha
$ws.Cells.Item(89, 10).Value = 2000          # J89 Volumen
$ws.Cells.Item(89, 11).Value = 800           # K89 Precio minimo
$ws.Cells.Item(89, 12).Value = 800           # L89 Precio maximo
$ws.Cells.Item(89, 13).Value = 800           # M89 Precio promedio ponderado
$ws.Cells.Item(89, 15).Value = "Región del Maule"   # O89 Origen
$ws.Cells.Item(89, 16).Value = 800           # P89 Precio $/Kg
